# Fill in actual game results (scores, overtime flag, win/loss, forecast
# accuracy) for rows 36-43 of Sheet1. These rows previously only had the
# scheduled matchup info (Date, Start, Away team, Home team, Arena); the
# games have since been played, so we add the final scores and derived
# columns, matching the pattern already used by the earlier rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Green fill color used to flag rows where the forecast favored the
# away (road) team and that forecast turned out correct. This matches
# the existing "FF00B050" fill already present in the workbook (style
# index 13 on K20/K21/K24/K27/K28/K32, etc.).
$upsetFillColor = 5287936   # RGB(0x00, 0xB0, 0x50)

$rows = @(
    @{ Row = 36; AwayPts = 106; HomePts = 113; Overtime = "NA";  Win = "Detroit Pistons";         Loss = "Charlotte Hornets";        Forecast = "Charlotte Hornets";        Correct = "No"  },
    @{ Row = 37; AwayPts = 118; HomePts = 107; Overtime = "NA";  Win = "Minnesota Timberwolves";   Loss = "Washington Wizards";       Forecast = "Washington Wizards";       Correct = "No"  },
    @{ Row = 38; AwayPts = 105; HomePts = 96;  Overtime = "NA";  Win = "Memphis Grizzlies";        Loss = "Miami Heat";               Forecast = "Memphis Grizzlies";        Correct = "Yes" },
    @{ Row = 39; AwayPts = 137; HomePts = 131; Overtime = "Yes"; Win = "Portland Trail Blazers";   Loss = "Houston Rockets";          Forecast = "Portland Trail Blazers";   Correct = "Yes" },
    @{ Row = 40; AwayPts = 116; HomePts = 126; Overtime = "NA";  Win = "Milwaukee Bucks";          Loss = "Cleveland Cavaliers";      Forecast = "Milwaukee Bucks";          Correct = "Yes" },
    @{ Row = 41; AwayPts = 132; HomePts = 109; Overtime = "NA";  Win = "Phoenix Suns";             Loss = "Dallas Mavericks";         Forecast = "Dallas Mavericks";         Correct = "No"  },
    @{ Row = 42; AwayPts = 140; HomePts = 114; Overtime = "NA";  Win = "Oklahoma City Thunder";    Loss = "San Antonio Spurs";        Forecast = "San Antonio Spurs";        Correct = "No"  },
    @{ Row = 43; AwayPts = 112; HomePts = 134; Overtime = "NA";  Win = "Golden State Warriors";    Loss = "Atlanta Hawks";            Forecast = "Atlanta Hawks";            Correct = "No"  }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # D: Away Pts, F: Home Pts (numeric, column style already applies #,##0 format)
    $ws.Cells.Item($rowNum, 4).Value = $r.AwayPts
    $ws.Cells.Item($rowNum, 6).Value = $r.HomePts

    # G: Overtime flag
    $ws.Cells.Item($rowNum, 7).Value = $r.Overtime

    # I: Win, J: Loss
    $ws.Cells.Item($rowNum, 9).Value = $r.Win
    $ws.Cells.Item($rowNum, 10).Value = $r.Loss

    # K: Forecasted winner
    $kCell = $ws.Cells.Item($rowNum, 11)
    $kCell.Value = $r.Forecast

    # Highlight the forecast cell green when the forecast was correct,
    # matching the existing styling convention used elsewhere in the sheet.
    if ($r.Correct -eq "Yes") {
        $kCell.Interior.Color = $upsetFillColor
    }

    # L: Correct (Yes/No)
    $ws.Cells.Item($rowNum, 12).Value = $r.Correct
}
